$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.661749362945557
$ws.Range("B1").Value = 2.899363994598389
$ws.Range("C1").Value = 3.6698899269104
$ws.Range("D1").Value = 1.425180673599243
$ws.Range("E1").Value = 0.9529051780700684
